# Weekly fruit/vegetable price log update:
# insert a new record as row 14, pushing the existing rows 14-34 down to 15-35.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(14).Insert()

$ws.Range("A14").Value = 1
$ws.Range("B14").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C14").Value = 'Arica y Parinacota'
$ws.Range("D14").Value = '2023-05-08'
$ws.Range("E14").Value = 15
$ws.Range("F14").Value = 100112013
$ws.Range("G14").Value = 'Alcachofa'
$ws.Range("H14").Value = 'Madrigal'
$ws.Range("I14").Value = 'Primera'
$ws.Range("J14").Value = 50
$ws.Range("K14").Value = 19000
$ws.Range("L14").Value = 20000
$ws.Range("M14").Value = 19600
$ws.Range("N14").Value = '$/caja 30 unidades'
$ws.Range("O14").Value = 'Región de Coquimbo'
$ws.Range("P14").Value = 653
$ws.Range("Q14").Value = 30
$ws.Range("R14").Value = 'Hortaliza'
